# Apply edit described in commit: extend the "Networks" TODO table (columns K:M)
# on the active sheet with three new rows of info, corresponding to new tasks
# about creating client avatars and changing when start/end positions are sent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: add a new task note in K16, with TODO status in M16 (matching
# existing style used elsewhere for "TODO" -> red fill).
$ws.Range("K16").Value = "Create client avatar"
$ws.Range("M16").Value = "TODO"
$ws.Range("M16").Interior.Color = 255

# Row 17: add new task + detail note in K17/L17, with TODO status in M17.
$ws.Range("K17").Value = "Change when new start/end positions sent"
$ws.Range("L17").Value = "Send on keypress"
$ws.Range("M17").Value = "TODO"
$ws.Range("M17").Interior.Color = 255

# Move the active selection to M17 (as recorded in the saved sheet view).
$ws.Range("M17").Select()
